$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("BCbVT-passenger", "BCbVT-freight")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # The existing G1 header ("nonroad vehicle") is being pushed out to a new
    # I1 column, G1 becomes "LPG vehicle", and a new H1 "hydrogen vehicle"
    # column is inserted in between. Capture the old text first.
    $nonroad = $ws.Range("G1").Value()

    $ws.Range("G1").Value = "LPG vehicle"

    # New H1 / I1 header cells - match the bold/right-aligned look of the
    # rest of the header row (B1:G1).
    $ws.Range("H1").Value = "hydrogen vehicle"
    $ws.Range("H1").Font.Bold = $true
    $ws.Range("H1").HorizontalAlignment = -4152

    $ws.Range("I1").Value = $nonroad
    $ws.Range("I1").Font.Bold = $true
    $ws.Range("I1").HorizontalAlignment = -4152

    # New leading label above the vehicle-type columns, word-wrapped.
    $ws.Range("A1").Value = "Battery Capacity (MW*hr/vehicle"
    $ws.Range("A1").WrapText = $true

    # Give the new header row extra height to fit the wrapped label.
    $ws.Rows.Item(1).RowHeight = 57

    # Size the two new columns similarly to the existing data columns.
    $ws.Columns.Item(8).ColumnWidth = 19.86328125
    $ws.Columns.Item(9).ColumnWidth = 18.265625

    # Fill the two new data columns (H, I) with zeros for rows 2-7, same as
    # the existing placeholder columns (C, D, E, G).
    for ($r = 2; $r -le 7; $r++) {
        $ws.Cells.Item($r, 8).Value = 0
        $ws.Cells.Item($r, 9).Value = 0
    }
}
